$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (e.g. "207.30", "0.0607").
# Coerce just those specific cells to Text format first so Excel stores them
# as the exact literal string instead of re-parsing/rounding them as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "25.304.29"
$ws.Range("E2").Value = "  -2.43%  "

$ws.Range("D3").Value = "1.566.18"
$ws.Range("E3").Value = "  -3.57%  "

$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").Value = "207.30"
$ws.Range("E5").Value = "  -2.98%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -4.95%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.0607"
$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.242"
$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").Value = "17.81"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").Value = "1.785.13"
$ws.Range("E12").Value = "  -3.47%  "

$ws.Range("D13").Value = "1.569.99"
$ws.Range("E13").Value = "  -3.33%  "

$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  -4.00%  "

$ws.Range("D15").Value = "0.506"
$ws.Range("E15").Value = "  -3.33%  "

$ws.Range("D16").Value = "25.311.71"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "59.43"
$ws.Range("E17").Value = "  -2.85%  "

$ws.Range("E18").Value = "  -2.92%  "

$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").Value = "185.38"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("D21").Value = "4.13"
$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  -2.83%  "

$ws.Range("E23").Value = "  -2.90%  "

$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("E25").Value = "  -4.16%  "

$ws.Range("D26").Value = "140.08"
$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("E27").Value = "  -6.85%  "

$ws.Range("E28").Value = "  -3.79%  "

$ws.Range("D29").Value = "14.84"
$ws.Range("E29").Value = "  -2.12%  "

$ws.Range("E30").Value = "  -6.03%  "

$ws.Range("D31").Value = "0.0464"
$ws.Range("E31").Value = "  -3.64%  "

$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  -3.54%  "

$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("E35").Value = "  -3.51%  "

$ws.Range("D36").Value = "1.090.98"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("E38").Value = "  -5.01%  "

$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.494"
$ws.Range("E40").Value = "  -4.79%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "0.770"
$ws.Range("E41").Value = "  -8.92%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.811"
$ws.Range("E42").Value = "  +5.72%  "

$ws.Range("D43").Value = "93.18"
$ws.Range("E43").Value = "  -4.84%  "

$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("D45").Value = "1.699.37"
$ws.Range("E45").Value = "  -3.46%  "

$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("D47").Value = "52.59"
$ws.Range("E47").Value = "  -3.39%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.44"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").Value = "  -4.63%  "

$ws.Range("E50").Value = "  -1.65%  "

$ws.Range("E51").Value = "  -0.50%  "

